$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels: shorten the voucher number column headers
$ws.Range("C4").Value = "CV NO."
$ws.Range("D4").Value = "JV NO."

# Adjust row height for the header row
$ws.Rows("4:4").RowHeight = 28.5

# Move the active selection
$ws.Range("C5").Select()
